$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.261.11"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.599.02"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'212.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.0606"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.244"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").Value = "'0.0856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "1.825.00"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "1.606.64"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "'4.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "'63.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "26.266.02"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'227.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'7.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.08%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'4.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "'8.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'145.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").Value = "'6.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").Value = "'15.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'0.0493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "1.447.60"
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("D34").Value = "'2.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").Value = "'0.0164"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("D43").Value = "'0.925"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "1.736.94"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'0.758"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("D46").Value = "'60.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").Value = "'87.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "'1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'0.0948"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
